$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 01:28"

# Update per-country COVID statistics (columns B,C,D,E,G,H) for the rows that changed

# Row 4
$ws.Range("B4").Value = 2775516
$ws.Range("C4").Value = 47663
$ws.Range("D4").Value = 1159577
$ws.Range("E4").Value = 1485176
$ws.Range("G4").Value = 641
$ws.Range("H4").Value = 130763

# Row 5
$ws.Range("B5").Value = 1453369
$ws.Range("C5").Value = 44884
$ws.Range("D5").Value = 826866
$ws.Range("E5").Value = 565790
$ws.Range("G5").Value = 1057
$ws.Range("H5").Value = 60713

# Row 28
$ws.Range("B28").Value = 67197
$ws.Range("C28").Value = 2667
$ws.Range("E28").Value = 42806
$ws.Range("G28").Value = 44
$ws.Range("H28").Value = 1351

# Row 31
$ws.Range("B31").Value = 58257
$ws.Range("C31").Value = 1825
$ws.Range("D31").Value = 27887
$ws.Range("E31").Value = 25794
$ws.Range("G31").Value = 49
$ws.Range("H31").Value = 4576

# Row 32
$ws.Range("B32").Value = 57770
$ws.Range("C32").Value = 1385
$ws.Range("D32").Value = 25595
$ws.Range("E32").Value = 29241
$ws.Range("G32").Value = 58
$ws.Range("H32").Value = 2934

# Row 41
$ws.Range("B41").Value = 41065
$ws.Range("C41").Value = 944
$ws.Range("D41").Value = 25319
$ws.Range("E41").Value = 15558
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 188

# Row 42
$ws.Range("B42").Value = 38511
$ws.Range("C42").Value = 997
$ws.Range("D42").Value = 10438
$ws.Range("E42").Value = 26803
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 1270

# Row 43
$ws.Range("B43").Value = 34775
$ws.Range("C43").Value = 382
$ws.Range("D43").Value = 21791
$ws.Range("E43").Value = 11507
$ws.Range("G43").Value = 14
$ws.Range("H43").Value = 1477

# Row 44
$ws.Range("B44").Value = 33550
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 15745
$ws.Range("E44").Value = 17174
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 631

# Row 45
$ws.Range("B45").Value = 33387
$ws.Range("C45").Value = 819
$ws.Range("D45").Value = 17904
$ws.Range("E45").Value = 14729
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 754

# Row 46
$ws.Range("B46").Value = 33219
$ws.Range("C46").Value = 1094
$ws.Range("D46").Value = 9340
$ws.Range("E46").Value = 22756
$ws.Range("G46").Value = 52
$ws.Range("H46").Value = 1123

# Row 47
$ws.Range("B47").Value = 31851
$ws.Range("C47").Value = 137
$ws.Range("D47").Value = 29200
$ws.Range("E47").Value = 686
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 1965

# Row 48
$ws.Range("B48").Value = 31836
$ws.Range("C48").Value = 319
$ws.Range("D48").Value = 15651
$ws.Range("E48").Value = 15411
$ws.Range("G48").Value = 28
$ws.Range("H48").Value = 774

# Row 49
$ws.Range("B49").Value = 27414
$ws.Range("C49").Value = 656
$ws.Range("D49").Value = 21948
$ws.Range("E49").Value = 5374
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 92

# Row 50
$ws.Range("B50").Value = 27296
$ws.Range("C50").Value = 326
$ws.Range("D50").Value = 19314
$ws.Range("E50").Value = 6315
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 1667

# Row 51
$ws.Range("B51").Value = 26484
$ws.Range("C51").Value = 790
$ws.Range("D51").Value = 10152
$ws.Range("E51").Value = 15729
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 603

# Row 52
$ws.Range("B52").Value = 26257
$ws.Range("C52").Value = 1013
$ws.Range("D52").Value = 17452
$ws.Range("E52").Value = 8483
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 322

# Row 53
$ws.Range("B53").Value = 26065
$ws.Range("C53").Value = 523
$ws.Range("D53").Value = 14563
$ws.Range("E53").Value = 11049
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 453

# Row 54
$ws.Range("B54").Value = 25477
$ws.Range("C54").Value = 4
$ws.Range("D54").Value = 23364
$ws.Range("E54").Value = 375
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 1738

# Row 56
$ws.Range("B56").Value = 18723
$ws.Range("C56").Value = 130
$ws.Range("D56").Value = 16731
$ws.Range("E56").Value = 1018
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 974

# Row 57
$ws.Range("B57").Value = 18134
$ws.Range("C57").Value = 393
$ws.Range("D57").Value = 13550
$ws.Range("E57").Value = 4467
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 117

# Row 58
$ws.Range("B58").Value = 18112
$ws.Range("C58").Value = 588
$ws.Range("D58").Value = 10061
$ws.Range("E58").Value = 7831
$ws.Range("G58").Value = 7
$ws.Range("H58").Value = 220

# Row 59
$ws.Range("B59").Value = 18096
$ws.Range("C59").Value = 687
$ws.Range("D59").Value = 3194
$ws.Range("E59").Value = 14129
$ws.Range("G59").Value = 27
$ws.Range("H59").Value = 773

# Row 60
$ws.Range("B60").Value = 17873
$ws.Range("C60").Value = 107
$ws.Range("D60").Value = 16491
$ws.Range("E60").Value = 677
$ws.Range("H60").Value = 705

# Row 69
$ws.Range("B69").Value = 12046
$ws.Range("C69").Value = 92
$ws.Range("E69").Value = 3900

# Row 87
$ws.Range("B87").Value = 5404
$ws.Range("C87").Value = 13
$ws.Range("D87").Value = 4346
$ws.Range("E87").Value = 1025
$ws.Range("H87").Value = 33

# Row 88
$ws.Range("B88").Value = 5394
$ws.Range("D88").Value = 2420
$ws.Range("E88").Value = 2932
$ws.Range("H88").Value = 42

# Row 96
$ws.Range("E96").Value = 2220
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 17

# Row 137
$ws.Range("B137").Value = 943
$ws.Range("C137").Value = 7
$ws.Range("D137").Value = 825
$ws.Range("E137").Value = 90
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 28

# Row 154
$ws.Range("B154").Value = 535
$ws.Range("C154").Value = 20
$ws.Range("D154").Value = 242
$ws.Range("E154").Value = 280
$ws.Range("H154").Value = 13

# Row 155
$ws.Range("B155").Value = 528
$ws.Range("D155").Value = 472
$ws.Range("E155").Value = 54
$ws.Range("H155").Value = 2

# Row 166
$ws.Range("B166").Value = 248
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 116
$ws.Range("E166").Value = 120
